# Weekly update of the Ají (pepper) vegetable price dataset for
# Mapocho Venta Directa de Santiago. Rows are re-populated with the
# latest reported values (Fecha, Variedad, Calidad, Volumen, Precios,
# Unidad de comercializacion, Precio $/Kg, Kg o Unidades).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace with data previously in row 5
$ws.Range("D2").Value = 44193
$ws.Range("H2").Value = 'Americana (o)'
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 46000
$ws.Range("L2").Value = 46000
$ws.Range("M2").Value = 46000
$ws.Range("N2").Value = '$/caja 15 kilos'
$ws.Range("P2").Value = 3067
$ws.Range("Q2").Value = 15

# Row 3: replace with data previously in row 8
$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 75000
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = 75000
$ws.Range("P3").Value = 3000

# Row 4: replace with data previously in row 9
$ws.Range("D4").Value = 44446
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 78000
$ws.Range("L4").Value = 78000
$ws.Range("M4").Value = 78000
$ws.Range("N4").Value = '$/caja 25 kilos'
$ws.Range("P4").Value = 3120
$ws.Range("Q4").Value = 25

# Row 5: replace with data previously in row 10
$ws.Range("D5").Value = 44446
$ws.Range("H5").Value = 'Inferno'
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 80000
$ws.Range("L5").Value = 80000
$ws.Range("M5").Value = 80000
$ws.Range("P5").Value = 5333

# Row 7: replace with data previously in row 17
$ws.Range("D7").Value = 44319
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 30000
$ws.Range("P7").Value = 1200

# Row 8: replace with data previously in row 11
$ws.Range("D8").Value = 44221
$ws.Range("J8").Value = 22
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 24545
$ws.Range("P8").Value = 982

# Row 9: replace with data previously in row 18
$ws.Range("D9").Value = 44343
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 36000
$ws.Range("L9").Value = 36000
$ws.Range("M9").Value = 36000
$ws.Range("P9").Value = 1440

# Row 10: replace with data previously in row 12
$ws.Range("D10").Value = 44421
$ws.Range("H10").Value = 'Americana (o)'
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 75000
$ws.Range("L10").Value = 75000
$ws.Range("M10").Value = 75000
$ws.Range("N10").Value = '$/caja 25 kilos'
$ws.Range("P10").Value = 3000
$ws.Range("Q10").Value = 25

# Row 11: replace with data previously in row 13
$ws.Range("D11").Value = 44326
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 30000
$ws.Range("P11").Value = 1200

# Row 12: replace with data previously in row 2
$ws.Range("D12").Value = 44553
$ws.Range("H12").Value = 'Inferno'
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 45000
$ws.Range("L12").Value = 45000
$ws.Range("M12").Value = 45000
$ws.Range("P12").Value = 1800

# Row 13: replace with data previously in row 7
$ws.Range("D13").Value = 44474
$ws.Range("J13").Value = 18
$ws.Range("K13").Value = 100000
$ws.Range("L13").Value = 100000
$ws.Range("M13").Value = 100000
$ws.Range("P13").Value = 4000

# Row 14: replace with data previously in row 16
$ws.Range("D14").Value = 44544
$ws.Range("H14").Value = 'Inferno'
$ws.Range("J14").Value = 12
$ws.Range("K14").Value = 35000
$ws.Range("L14").Value = 35000
$ws.Range("M14").Value = 35000
$ws.Range("P14").Value = 1400

# Row 16: replace with data previously in row 3
$ws.Range("D16").Value = 44449
$ws.Range("H16").Value = 'Americana (o)'
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 80000
$ws.Range("L16").Value = 80000
$ws.Range("M16").Value = 80000
$ws.Range("P16").Value = 3200

# Row 17: replace with data previously in row 4
$ws.Range("D17").Value = 44449
$ws.Range("I17").Value = 'Segunda'
$ws.Range("K17").Value = 75000
$ws.Range("L17").Value = 75000
$ws.Range("M17").Value = 75000
$ws.Range("N17").Value = '$/caja 15 kilos'
$ws.Range("P17").Value = 5000
$ws.Range("Q17").Value = 15

# Row 18: replace with data previously in row 14
$ws.Range("D18").Value = 44460
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 95000
$ws.Range("L18").Value = 95000
$ws.Range("M18").Value = 95000
$ws.Range("P18").Value = 3800
